$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 22.97
$ws.Range("P2").Value = 40.86
$ws.Range("Q2").Value = 63.83
$ws.Range("O3").Value = 23.13
$ws.Range("P3").Value = 40.47
$ws.Range("Q3").Value = 63.6
$ws.Range("O4").Value = 22.89
$ws.Range("P4").Value = 41.8
$ws.Range("Q4").Value = 64.69
$ws.Range("O5").Value = 24.85
$ws.Range("P5").Value = 39.65
$ws.Range("Q5").Value = 64.5
$ws.Range("O6").Value = 23.93
$ws.Range("P6").Value = 41.26
$ws.Range("Q6").Value = 65.19
$ws.Range("O7").Value = 23.88
$ws.Range("P7").Value = 40.11
$ws.Range("Q7").Value = 63.99
$ws.Range("O8").Value = 22.9
$ws.Range("P8").Value = 41.44
$ws.Range("Q8").Value = 64.34
$ws.Range("O9").Value = 23.26
$ws.Range("P9").Value = 39.95
$ws.Range("Q9").Value = 63.21
$ws.Range("O10").Value = 23.57
$ws.Range("P10").Value = 39.38
$ws.Range("Q10").Value = 62.95
$ws.Range("O11").Value = 23.72
$ws.Range("P11").Value = 40.56
$ws.Range("Q11").Value = 64.28
$ws.Range("O12").Value = 23.7
$ws.Range("P12").Value = 41.45
$ws.Range("Q12").Value = 65.15
$ws.Range("O13").Value = 22.69
$ws.Range("P13").Value = 42.22
$ws.Range("Q13").Value = 64.91
$ws.Range("O14").Value = 23.76
$ws.Range("P14").Value = 43.82
$ws.Range("Q14").Value = 67.58
$ws.Range("O15").Value = 23.51
$ws.Range("P15").Value = 43.03
$ws.Range("Q15").Value = 66.54
$ws.Range("O16").Value = 22.98
$ws.Range("P16").Value = 43.29
$ws.Range("Q16").Value = 66.27
$ws.Range("O17").Value = 25.65
$ws.Range("P17").Value = 40.12
$ws.Range("Q17").Value = 65.77
$ws.Range("O18").Value = 24.95
$ws.Range("P18").Value = 40.24
$ws.Range("Q18").Value = 65.19
$ws.Range("O19").Value = 24
$ws.Range("P19").Value = 41.74
$ws.Range("Q19").Value = 65.74
$ws.Range("O20").Value = 23.88
$ws.Range("P20").Value = 40.8
$ws.Range("Q20").Value = 64.68
$ws.Range("O21").Value = 23.05
$ws.Range("P21").Value = 42.04
$ws.Range("Q21").Value = 65.09
$ws.Range("O22").Value = 23.05
$ws.Range("P22").Value = 41.96
$ws.Range("Q22").Value = 65.01
$ws.Range("O23").Value = 22.56
$ws.Range("P23").Value = 42.25
$ws.Range("Q23").Value = 64.81
$ws.Range("O24").Value = 22.67
$ws.Range("P24").Value = 42.47
$ws.Range("Q24").Value = 65.14
$ws.Range("O25").Value = 21.94
$ws.Range("P25").Value = 42.39
$ws.Range("Q25").Value = 64.33
$ws.Range("O26").Value = 21.1
$ws.Range("P26").Value = 44.72
$ws.Range("Q26").Value = 65.82
$ws.Range("O27").Value = 22.05
$ws.Range("P27").Value = 43.96
$ws.Range("Q27").Value = 66.01
$ws.Range("O28").Value = 21.7
$ws.Range("P28").Value = 43.61
$ws.Range("Q28").Value = 65.31
$ws.Range("O29").Value = 22.48
$ws.Range("P29").Value = 40.88
$ws.Range("Q29").Value = 63.36
$ws.Range("O30").Value = 22.57
$ws.Range("P30").Value = 40.82
$ws.Range("Q30").Value = 63.39
$ws.Range("O31").Value = 22.5
$ws.Range("P31").Value = 40.65
$ws.Range("Q31").Value = 63.15
